# Generate Report for Handoff
# Rotates the localized file's identity from the "92275e07-...".md GUID to
# the new "c99586a5-...".md GUID, refreshes the zh-cn / de-de handoff xlf
# names + timestamps, and clears out the "already handed back" columns for
# both language sheets now that a brand new handoff round has started.

$wb = $excel.ActiveWorkbook

$oldGuid = "92275e07-faf6-479a-a38a-950c2959146a"
$newGuid = "c99586a5-0ad5-41ed-b386-4315f52cd227"

$oldHash = "4390e6f653466e5aead1aea3810d6008917612d5"
$newHash = "82927371a8bc552d1dc93c50a06976ba1a0a5dbc"

$overviewAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16e0c8d940025351255dde198b7dc0af6c34f399/e2e/$oldGuid.md"
$zhcnAddr     = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d2286968565ee492812dfe45fa996b0ea6d4b2dd/e2e/$oldGuid.md"
$dedeAddr     = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dc4edcbb49b2b2b9a32dbd9170bf6bfb4c4edf3f/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-04 07:06:09"

# Refresh B2's hyperlink display text, keeping its target address intact.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewAddr, "", "", "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-04 07:06:00"
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# The file has not been handed back yet under its new identity, so the
# "Latest Target File" / "Latest Handback File" columns are cleared and
# lose their hyperlink + hyperlink styling.
$wsZhCn.Range("I2").Hyperlinks.Delete()
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""

# Re-create the remaining A2 hyperlink (display text only changes; address
# is unchanged) since deleting hyperlinks above clears the whole sheet.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnAddr, "", "", "$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-04 07:06:09"
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("I2").Hyperlinks.Delete()
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeAddr, "", "", "$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# Column I / J on the language sheets shrink now that the long file-name
# hyperlinks are gone; auto-fit them to the remaining (shorter) content.
# ---------------------------------------------------------------------
$wsZhCn.Columns.Item(9).AutoFit() | Out-Null
$wsZhCn.Columns.Item(10).AutoFit() | Out-Null
$wsDeDe.Columns.Item(9).AutoFit() | Out-Null
$wsDeDe.Columns.Item(10).AutoFit() | Out-Null
